$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Gestion du temps" update: the descriptions in column C for rules
# RG04 (row 62), RG05 (row 63) and RG06 (row 64) were missing; they had
# been entered three rows too low (against RG07/RG08/RG09 on rows
# 65-67). Move each description up to line up with the correct rule.
$ws.Cells.Item(62, 3).Value = $ws.Cells.Item(65, 3).Value2
$ws.Cells.Item(63, 3).Value = $ws.Cells.Item(66, 3).Value2
$ws.Cells.Item(64, 3).Value = $ws.Cells.Item(67, 3).Value2

# Clear the old (now duplicated / misplaced) values.
$ws.Cells.Item(65, 3).ClearContents()
$ws.Cells.Item(66, 3).ClearContents()
$ws.Cells.Item(67, 3).ClearContents()

# Leave the selection where the user ended up working, on C65.
[void]$ws.Range("C65").Select()
